# Update column G ("K") values on the active sheet (Sheet1)
# per the diff: rows 2-11, values change as follows:
#   G2: 20 -> 5
#   G3: 5  -> 0
#   G4: 18 -> 6
#   G5: 22 -> 5
#   G6: 26 -> 4
#   G7: 16 -> 2
#   G8: 17 -> 5
#   G9: 10 -> 3
#   G10: 9 -> 3
#   G11: 4 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 0
    4  = 6
    5  = 5
    6  = 4
    7  = 2
    8  = 5
    9  = 3
    10 = 3
    11 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
